# Adding more experiments (hidden_dim = 100, 300, 600 runs) to the
# "experiments_colors" results sheet, and tweaking the description of the
# experiment 8 row.
#
# Final row order/content (rows 2-15, columns B..E) mirrors the target
# workbook exactly; row 3 is the edited "experiment 8" row and rows 2, 5,
# 7 and 15 are the newly added experiments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LSTM 1 layer + hidden_dim=100 (new)
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 0.809962568384681
$ws.Range("D2").Value = 0.632187527409257
$ws.Range("E2").Value = "LSTM 1 layer + hidden_dim=100"

# Row 3 - experiment 8, description updated
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 0.808810826374892
$ws.Range("D3").Value = 0.633431206619349
$ws.Range("E3").Value = "Based on experiment 1 + ColorizedInputDescriber (targets added to embeddings)"

# Row 4 - experiment 1 (unchanged, shifted down)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.803627987330837
$ws.Range("D4").Value = 0.630708944716922
$ws.Range("E4").Value = "Reusing GloVe embedding + Bert tokenization "

# Row 5 - RNN 3 layers encoder/decoder (new)
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 0.796717535272099
$ws.Range("D5").Value = 0.629499048876544
$ws.Range("E5").Value = "RNN 3 layers encoder/decoder"

# Row 6 - experiment 3 (unchanged, shifted down)
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 0.795565793262309
$ws.Range("D6").Value = 0.659179779992603
$ws.Range("E6").Value = "Glove + twitter tokenization"

# Row 7 - LSTM 1 layer (new)
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 0.795277857759861
$ws.Range("D7").Value = 0.633255077309909
$ws.Range("E7").Value = "LSTM 1 layer"

# Row 8 - experiment 5 (unchanged, shifted down)
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 0.790095018715807
$ws.Range("D8").Value = 0.629040307809219
$ws.Range("E8").Value = "bert+bert tokens"

# Row 9 - experiment 6 (unchanged, shifted down, no description)
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 0.778289663115462
$ws.Range("D9").Value = 0.657156478324349
$ws.Range("E9").Value = ""

# Row 10 - experiment 0.1 (unchanged, shifted down)
$ws.Range("B10").Value = 0.1
$ws.Range("C10").Value = 0.764180823495537
$ws.Range("D10").Value = 0.677939279248182
$ws.Range("E10").Value = "Baseline with dev_mod based on whitspace tokenization and colors represented by fourier transform"

# Row 11 - experiment 2 (unchanged, shifted down)
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 0.757558306939245
$ws.Range("D11").Value = 0.6221401283633
$ws.Range("E11").Value = "Glove + white-space tokenization"

# Row 12 - experiment 4 (unchanged, shifted down, no description)
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 0.689893463864094
$ws.Range("D12").Value = 0.707073602479846
$ws.Range("E12").Value = ""

# Row 13 - experiment 7 (unchanged, shifted down, no description)
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 0.689317592859199
$ws.Range("D13").Value = 0.706977637009309
$ws.Range("E13").Value = ""

# Row 14 - experiment 0 (unchanged, shifted down)
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.399654477397063
$ws.Range("D14").Value = 0.057112006910452
$ws.Range("E14").Value = "Baseline with simple dev_mod and colors not represented by fourier transform"

# Row 15 - LSTM 1 layer + hidden_dim=300 (new)
$ws.Range("B15").Value = 11
$ws.Range("C15").Value = 0.813417794414051
$ws.Range("D15").Value = 0.631379150353065
$ws.Range("E15").Value = "LSTM 1 layer + hidden_dim=300"

# Column D grew a little wider (no longer auto "best fit") to comfortably
# show the longer hidden_dim labels.
$ws.Columns.Item(4).ColumnWidth = 17.5

# Move the active selection down past the new rows, matching where the
# cursor ended up after the edits.
$ws.Range("E17").Select()
